$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Update the Quantity value; dependent formulas in D27:D30 reference $C$23
# and will recalculate automatically.
$ws.Range("C23").Value = 8

$excel.CalculateFullRebuild()
